# "spirit flow spells finished"
# Update the spirit-flow spell modifier selections on Sheet2:
#  - C3: Barrier -> Spirit Control
#  - D3: All -> Drain
#  - E3: Roll Count -> None
#  - F3: Channel -> None
#  - B5: Barrier -> Spirit Recover
# The dependent INDEX/MATCH + CONCAT formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("C3").Value = "Spirit Control"
$ws.Range("D3").Value = "Drain"
$ws.Range("E3").Value = "None"
$ws.Range("F3").Value = "None"
$ws.Range("B5").Value = "Spirit Recover"

$excel.Calculate()

# Reflect the user's last selected cell on the sheet.
$ws.Activate()
$ws.Range("F6").Select()
